$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the SPEC min/max values (K20:L29) while keeping their number formatting.
$ws.Range("K20:L29").ClearContents()

# Move the active selection from N11 to H4, matching the saved workbook state.
$ws.Range("H4").Select()
